# Apply the "Drop 3" test-data changes to donationInfo.xlsx
#
# Summary of the change (per the commit diff):
#  - donations (sheet1): just a cursor-position change (C30 -> E3)
#  - donationShortDetail (sheet2): widen column D, move the selection, and
#    replace the "UNT01 / UNT01-2020122" request/donation id pair with a
#    new "SDK11006 / SDK11006-2020190" pair in row 3
#  - orderTests (sheet4): move the selection, replace the same pair in row 3
#    with the new SDK11006 pair (row 8 keeps its existing UNT01 pair), and
#    stop being the active/selected tab
#  - orderTestsStatusHistory (sheet5): move the selection, replace the
#    "ABV9005 / ABV9005-2020121" pair in row 3 with the new SDK11006 pair,
#    and become the active/selected tab

$wb = $excel.ActiveWorkbook

# --- donations -----------------------------------------------------------
$donations = $wb.Worksheets.Item("donations")
$donations.Range("E3").Select() | Out-Null

# --- donationShortDetail --------------------------------------------------
$shortDetail = $wb.Worksheets.Item("donationShortDetail")
$shortDetail.Range("A3").Value = "SDK11006"
$shortDetail.Range("B3").Value = "SDK11006-2020190"
$shortDetail.Columns.Item(4).ColumnWidth = 24.7
$shortDetail.Range("A3:B3").Select() | Out-Null

# --- orderTests ------------------------------------------------------------
$orderTests = $wb.Worksheets.Item("orderTests")
$orderTests.Range("A3").Value = "SDK11006"
$orderTests.Range("B3").Value = "SDK11006-2020190"
$orderTests.Range("F19").Select() | Out-Null

# --- orderTestsStatusHistory -----------------------------------------------
$statusHistory = $wb.Worksheets.Item("orderTestsStatusHistory")
$statusHistory.Range("A3").Value = "SDK11006"
$statusHistory.Range("B3").Value = "SDK11006-2020190"
$statusHistory.Range("A3:B3").Select() | Out-Null

# orderTestsStatusHistory becomes the active/selected sheet (was orderTests)
$statusHistory.Activate() | Out-Null
